# Scheduled runner update: refresh computed profit figures (columns H-N)
# across the per-job "Profits" sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 2450
$ws.Cells.Item(29, 10).Value = 2645.4546
$ws.Cells.Item(29, 12).Value = 7936.3638
$ws.Cells.Item(29, 14).Value = -8498.363799999999

$ws.Cells.Item(38, 8).Value = 633.7778
$ws.Cells.Item(38, 9).Value = 350.66666
$ws.Cells.Item(38, 10).Value = 1200
$ws.Cells.Item(38, 11).Value = 1051.99998
$ws.Cells.Item(38, 12).Value = 3600
$ws.Cells.Item(38, 13).Value = -679.9999800000001
$ws.Cells.Item(38, 14).Value = -4344

$ws.Cells.Item(58, 8).Value = 1173
$ws.Cells.Item(58, 10).Value = 4500
$ws.Cells.Item(58, 12).Value = 13500
$ws.Cells.Item(58, 14).Value = -13800

$ws.Cells.Item(80, 8).Value = 1044.5
$ws.Cells.Item(80, 9).Value = 2030.6666
$ws.Cells.Item(80, 10).Value = 715.7778
$ws.Cells.Item(80, 11).Value = 6091.9998
$ws.Cells.Item(80, 12).Value = 2147.3334
$ws.Cells.Item(80, 13).Value = -5093.9998
$ws.Cells.Item(80, 14).Value = -4143.3334

$ws.Cells.Item(83, 8).Value = 1044.5
$ws.Cells.Item(83, 9).Value = 2030.6666
$ws.Cells.Item(83, 10).Value = 715.7778
$ws.Cells.Item(83, 11).Value = 18275.9994
$ws.Cells.Item(83, 12).Value = 6442.000199999999
$ws.Cells.Item(83, 13).Value = -13283.9994
$ws.Cells.Item(83, 14).Value = -16426.0002

$ws.Cells.Item(96, 8).Value = 629.4286
$ws.Cells.Item(96, 9).Value = 551.5
$ws.Cells.Item(96, 10).Value = 733.3333
$ws.Cells.Item(96, 11).Value = 1654.5
$ws.Cells.Item(96, 12).Value = 2199.9999
$ws.Cells.Item(96, 13).Value = -281.5
$ws.Cells.Item(96, 14).Value = -4945.9999

$ws.Cells.Item(132, 8).Value = 15161374
$ws.Cells.Item(132, 9).Value = 22231322
$ws.Cells.Item(132, 10).Value = 11487.429
$ws.Cells.Item(132, 11).Value = 66693966
$ws.Cells.Item(132, 12).Value = 34462.287
$ws.Cells.Item(132, 13).Value = -66691436
$ws.Cells.Item(132, 14).Value = -39522.287

$ws.Cells.Item(137, 8).Value = 1635.3334
$ws.Cells.Item(137, 9).Value = 1001.7143
$ws.Cells.Item(137, 10).Value = 2038.5454
$ws.Cells.Item(137, 11).Value = 3005.1429
$ws.Cells.Item(137, 12).Value = 6115.6362
$ws.Cells.Item(137, 13).Value = -455.1428999999998
$ws.Cells.Item(137, 14).Value = -11215.6362

$ws.Cells.Item(138, 8).Value = 1427.3469
$ws.Cells.Item(138, 9).Value = 628.5
$ws.Cells.Item(138, 10).Value = 1978.2759
$ws.Cells.Item(138, 11).Value = 1885.5
$ws.Cells.Item(138, 12).Value = 5934.8277
$ws.Cells.Item(138, 13).Value = 3254.5
$ws.Cells.Item(138, 14).Value = -16214.8277

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2903.65
$ws.Cells.Item(32, 9).Value = 2636.9856
$ws.Cells.Item(32, 11).Value = 2636.9856
$ws.Cells.Item(32, 13).Value = -2349.9856

$ws.Cells.Item(45, 8).Value = 1137.7646
$ws.Cells.Item(45, 9).Value = 1103.2307
$ws.Cells.Item(45, 10).Value = 1250
$ws.Cells.Item(45, 11).Value = 1103.2307
$ws.Cells.Item(45, 12).Value = 1250
$ws.Cells.Item(45, 13).Value = -726.2307000000001
$ws.Cells.Item(45, 14).Value = -2004

$ws.Cells.Item(61, 8).Value = 1026.6444
$ws.Cells.Item(61, 9).Value = 863.3077
$ws.Cells.Item(61, 11).Value = 863.3077
$ws.Cells.Item(61, 13).Value = -651.3077

$ws.Cells.Item(132, 8).Value = 1466.6875
$ws.Cells.Item(132, 9).Value = 1176.2858
$ws.Cells.Item(132, 11).Value = 3528.8574
$ws.Cells.Item(132, 13).Value = -998.8574000000003

$ws.Cells.Item(136, 8).Value = 1026.6444
$ws.Cells.Item(136, 9).Value = 863.3077
$ws.Cells.Item(136, 11).Value = 2589.9231
$ws.Cells.Item(136, 13).Value = -39.92309999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(110, 8).Value = 42913.4
$ws.Cells.Item(110, 10).Value = 42913.4
$ws.Cells.Item(110, 12).Value = 42913.4
$ws.Cells.Item(110, 14).Value = -51093.4

$ws.Cells.Item(134, 8).Value = 4939.4375
$ws.Cells.Item(134, 9).Value = 1196.5652
$ws.Cells.Item(134, 10).Value = 14504.556
$ws.Cells.Item(134, 11).Value = 3589.6956
$ws.Cells.Item(134, 12).Value = 43513.66800000001
$ws.Cells.Item(134, 13).Value = -1054.6956
$ws.Cells.Item(134, 14).Value = -48583.66800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(33, 8).Value = 1100
$ws.Cells.Item(33, 9).Value = 1100
$ws.Cells.Item(33, 11).Value = 1100
$ws.Cells.Item(33, 13).Value = -721

$ws.Cells.Item(106, 8).Value = 50000
$ws.Cells.Item(106, 10).Value = 50000
$ws.Cells.Item(106, 12).Value = 50000
$ws.Cells.Item(106, 14).Value = -52524

$ws.Cells.Item(134, 8).Value = 1071.35
$ws.Cells.Item(134, 9).Value = 760.5
$ws.Cells.Item(134, 10).Value = 1796.6666
$ws.Cells.Item(134, 11).Value = 2281.5
$ws.Cells.Item(134, 12).Value = 5389.9998
$ws.Cells.Item(134, 13).Value = 253.5
$ws.Cells.Item(134, 14).Value = -10459.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 2143.7917
$ws.Cells.Item(39, 10).Value = 1907.1904
$ws.Cells.Item(39, 12).Value = 5721.5712
$ws.Cells.Item(39, 14).Value = -6309.5712

$ws.Cells.Item(55, 8).Value = 2500
$ws.Cells.Item(55, 10).Value = 2875
$ws.Cells.Item(55, 12).Value = 8625
$ws.Cells.Item(55, 14).Value = -8979

$ws.Cells.Item(136, 8).Value = 1727.625
$ws.Cells.Item(136, 9).Value = 1091.8182
$ws.Cells.Item(136, 10).Value = 3126.4
$ws.Cells.Item(136, 11).Value = 3275.4546
$ws.Cells.Item(136, 12).Value = 9379.200000000001
$ws.Cells.Item(136, 13).Value = 1824.5454
$ws.Cells.Item(136, 14).Value = -19579.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 173.11111
$ws.Cells.Item(2, 9).Value = 176.33333
$ws.Cells.Item(2, 11).Value = 176.33333
$ws.Cells.Item(2, 13).Value = -63.33332999999999

$ws.Cells.Item(97, 8).Value = 1010
$ws.Cells.Item(97, 9).Value = 1010
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 1010
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).Value = -514
$ws.Cells.Item(97, 14).ClearContents()

$ws.Cells.Item(122, 8).Value = 1463.4166
$ws.Cells.Item(122, 9).Value = 1446.7368
$ws.Cells.Item(122, 11).Value = 4340.2104
$ws.Cells.Item(122, 13).Value = -1890.2104

$ws.Cells.Item(132, 8).Value = 2303.75
$ws.Cells.Item(132, 9).Value = 1729.4286
$ws.Cells.Item(132, 10).Value = 2750.4443
$ws.Cells.Item(132, 11).Value = 5188.2858
$ws.Cells.Item(132, 12).Value = 8251.332900000001
$ws.Cells.Item(132, 13).Value = -2658.2858
$ws.Cells.Item(132, 14).Value = -13311.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2635.75
$ws.Cells.Item(40, 9).Value = 1634.2727
$ws.Cells.Item(40, 11).Value = 1634.2727
$ws.Cells.Item(40, 13).Value = -1498.2727

$ws.Cells.Item(63, 8).Value = 20000
$ws.Cells.Item(63, 10).Value = 20000
$ws.Cells.Item(63, 12).Value = 20000
$ws.Cells.Item(63, 14).Value = -21498

$ws.Cells.Item(66, 8).Value = 20000
$ws.Cells.Item(66, 10).Value = 20000
$ws.Cells.Item(66, 12).Value = 60000
$ws.Cells.Item(66, 14).Value = -67488

$ws.Cells.Item(132, 8).Value = 21045.941
$ws.Cells.Item(132, 9).Value = 944
$ws.Cells.Item(132, 10).Value = 65018.938
$ws.Cells.Item(132, 11).Value = 2832
$ws.Cells.Item(132, 12).Value = 195056.814
$ws.Cells.Item(132, 13).Value = -302
$ws.Cells.Item(132, 14).Value = -200116.814

$ws.Cells.Item(136, 8).Value = 1666.1666
$ws.Cells.Item(136, 9).Value = 1569.1428
$ws.Cells.Item(136, 10).Value = 1802
$ws.Cells.Item(136, 11).Value = 4707.428400000001
$ws.Cells.Item(136, 12).Value = 5406
$ws.Cells.Item(136, 13).Value = -2157.428400000001
$ws.Cells.Item(136, 14).Value = -10506

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1377.7727
$ws.Cells.Item(132, 9).Value = 841.86664
$ws.Cells.Item(132, 10).Value = 2526.1428
$ws.Cells.Item(132, 11).Value = 2525.59992
$ws.Cells.Item(132, 12).Value = 7578.428400000001
$ws.Cells.Item(132, 13).Value = 4.400080000000344
$ws.Cells.Item(132, 14).Value = -12638.4284
